$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Données")

# Add new "Obligatoire" header in column D, matching the header style of A1:C1
$ws.Cells.Item(1, 4).Value = "Obligatoire"
$ws.Cells.Item(1, 1).Copy()
$ws.Cells.Item(1, 4).PasteSpecial(-4122)

# Fill column D ("Obligatoire") with "Oui" for every data row, matching the
# row's existing data style (same as column A/B/C for that row)
$lastRow = 28
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = "Oui"
    $ws.Cells.Item($r, 1).Copy()
    $ws.Cells.Item($r, 4).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
